$wb = $excel.ActiveWorkbook

# --- Sheet: species_in_model ---
$ws1 = $wb.Worksheets.Item("species_in_model")

# Row 2 (Rremin / Porg / reactant)
$ws1.Range("F2").Value = "Porg"
$ws1.Range("J2").Value = 1

# Row 3 (Rremin / H3PO4 -> TH3PO4 / product)
$ws1.Range("B3").Value = "TH3PO4"
$ws1.Range("F3").Value = "TH3PO4,P_ads,TH3PO4_ads,TH3PO4_ads,TH3PO4_dis"
$ws1.Range("J3").Value = 2

# Row 4 (Rpre / H3PO4 -> TH3PO4 / reactant)
$ws1.Range("B4").Value = "TH3PO4"
$ws1.Range("F4").Value = "TH3PO4,P_ads,TH3PO4_ads,TH3PO4_ads,TH3PO4_dis"
$ws1.Range("J4").Value = 2

# --- Sheet: transport_parameters ---
$ws5 = $wb.Worksheets.Item("transport_parameters")

# Swap row2 / row3 parameter names (TH3PO4ID <-> PorgID)
$ws5.Range("A2").Value = "PorgID"
$ws5.Range("A3").Value = "TH3PO4ID"

# Rename transport-matrix parameters
$ws5.Range("A5").Value = "AmTH3PO4_dis"
$ws5.Range("A6").Value = "AmTH3PO4_ads"

# Rename boundary-condition parameters
$ws5.Range("A10").Value = "BcAmTH3PO4_dis"
$ws5.Range("A11").Value = "BcCmTH3PO4_dis"
$ws5.Range("A12").Value = "BcAmTH3PO4_ads"
$ws5.Range("A13").Value = "BcCmTH3PO4_ads"

# Rename bioirrigation parameter
$ws5.Range("A15").Value = "TH3PO4_dis0"

# Row 16 becomes "dstopw" / "adsorption" (was row 17); drop the old
# K_ads / adsorption row (row 16) by shifting row 17 up and deleting row 17.
$ws5.Range("A16").Value = "dstopw"
$ws5.Range("B16").Value = "adsorption"
$ws5.Rows.Item(17).Delete()

# --- Sheet: reaction_parameters ---
$ws6 = $wb.Worksheets.Item("reaction_parameters")

# K_ads row: adsorption -> speciation, comment updated
$ws6.Range("B2").Value = "speciation"
$ws6.Range("C2").Value = "TH3PO4_dis,P_ads"

# dstopw row: adsorption -> speciation, comment updated
$ws6.Range("B3").Value = "speciation"
$ws6.Range("C3").Value = "TH3PO4_dis"
